# Generate Report for Handback
# Updates the localization-status workbook to reflect that the
# 4a452103-27ad-4a63-a2bc-f2f24f1b01bc.md file has now been handed back
# (in sync with en-US) for both the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"

# Overview sheet: row 3 corresponds to 4a452103-27ad-4a63-a2bc-f2f24f1b01bc.md
# Update the zh-cn (E3) and de-de (F3) status columns.
$overview.Range("E3").Value = $handedBack
$overview.Range("F3").Value = $handedBack

# zh-cn detail sheet: row 3 is the 4a452103... file.
# Status -> Handed back, Latest Handback DateTime -> new timestamp,
# Error Detail -> cleared (kept as an empty-string cell, not removed).
$zhcn.Range("C3").Value = $handedBack
$zhcn.Range("K3").Value = "2016-08-30 15:03:43"
$zhcn.Range("P3").Formula = "'"
$zhcn.Range("P3").Style = "Normal"

# de-de detail sheet: row 3 is the 4a452103... file.
# Status -> Handed back, Latest Handback DateTime -> new timestamp,
# Error Detail -> cleared (kept as an empty-string cell, not removed).
$dede.Range("C3").Value = $handedBack
$dede.Range("K3").Value = "2016-08-30 15:03:50"
$dede.Range("P3").Formula = "'"
$dede.Range("P3").Style = "Normal"
